$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("P2").Value = 972
$ws.Range("D3").Value = 44637
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 861
$ws.Range("D4").Value = 44819
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 15000
$ws.Range("P4").Value = 833
$ws.Range("D5").Value = 44630
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15500
$ws.Range("P5").Value = 861
$ws.Range("D6").Value = 44635
$ws.Range("J6").Value = 100
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 861
$ws.Range("D7").Value = 44799
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 861
$ws.Range("D8").Value = 44782
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("P8").Value = 972
$ws.Range("D9").Value = 44754
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 16500
$ws.Range("P9").Value = 917
$ws.Range("D10").Value = 44645
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("P10").Value = 861
$ws.Range("D11").Value = 44775
$ws.Range("D12").Value = 44658
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 861
$ws.Range("D13").Value = 44791
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17500
$ws.Range("P13").Value = 972
$ws.Range("D14").Value = 44651
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 861
$ws.Range("D15").Value = 44642
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 15500
$ws.Range("P15").Value = 861
$ws.Range("D16").Value = 44790
$ws.Range("J16").Value = 60
$ws.Range("D17").Value = 44804
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 861
$ws.Range("D18").Value = 44832
$ws.Range("D19").Value = 44763
$ws.Range("J19").Value = 80
$ws.Range("D20").Value = 44771
$ws.Range("D21").Value = 44785
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("P21").Value = 972
$ws.Range("D22").Value = 44761
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("P22").Value = 972
$ws.Range("D23").Value = 44811
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("P23").Value = 806
$ws.Range("D24").Value = 44809
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("P24").Value = 806
$ws.Range("D25").Value = 44664
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15500
$ws.Range("P25").Value = 861
$ws.Range("D26").Value = 44628
$ws.Range("J26").Value = 60
$ws.Range("D27").Value = 44847
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 17000
$ws.Range("P27").Value = 944
$ws.Range("D28").Value = 44656
$ws.Range("J28").Value = 100
$ws.Range("D29").Value = 44830
$ws.Range("J29").Value = 60
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 17000
$ws.Range("P29").Value = 944
$ws.Range("D30").Value = 44659
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 80
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("P30").Value = 861
$ws.Range("D31").Value = 44813
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 14000
$ws.Range("M31").Value = 14500
$ws.Range("P31").Value = 806
$ws.Range("D32").Value = 44769
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17500
$ws.Range("P32").Value = 972
$ws.Range("D33").Value = 44797
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 16500
$ws.Range("P33").Value = 917
$ws.Range("D34").Value = 44818
$ws.Range("J34").Value = 60
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("P34").Value = 833
$ws.Range("D35").Value = 44384
$ws.Range("J35").Value = 120
$ws.Range("K35").Value = 17000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 17500
$ws.Range("P35").Value = 972
$ws.Range("D36").Value = 44384
$ws.Range("I36").Value = "Segunda"
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = 15000
$ws.Range("P36").Value = 833
$ws.Range("D37").Value = 44649
$ws.Range("J37").Value = 60
